$d = $word.ActiveDocument

# ---- Page margins (APA format: 1 inch all around) ----
$ps = $d.PageSetup
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.LeftMargin = 72
$ps.RightMargin = 72
$ps.HeaderDistance = 35.3
$ps.FooterDistance = 35.3

# ---- Normal style (APA formatting: Times New Roman 12pt, double spaced,
#      0.5" first line indent) ----
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.LineSpacingRule = 2
$normal.ParagraphFormat.FirstLineIndent = 36
$normal.Font.Name = "Times New Roman"
$normal.Font.Size = 12

# ---- New "heading 1" / Ttulo1 paragraph style + its linked character style ----
$ttulo1 = $d.Styles.Add("Ttulo1", 1)
$ttulo1Car = $d.Styles.Add("Ttulo1Car", 2)

$ttulo1.NameLocal = "heading 1"
$ttulo1.BaseStyle = $normal
$ttulo1.NextParagraphStyle = $normal
$ttulo1.LinkStyle = $ttulo1Car
$ttulo1.Priority = 9
$ttulo1.QuickStyle = $true

$ttulo1.ParagraphFormat.KeepWithNext = $true
$ttulo1.ParagraphFormat.KeepTogether = $true
$ttulo1.ParagraphFormat.SpaceBefore = 12
$ttulo1.ParagraphFormat.SpaceAfter = 0
$ttulo1.ParagraphFormat.FirstLineIndent = 0
$ttulo1.ParagraphFormat.OutlineLevel = 1

$ttulo1.Font.Bold = $true
$ttulo1.Font.Size = 14
$ttulo1.Font.SizeBi = 16

$ttulo1Car.NameLocal = "Título 1 Car"
$ttulo1Car.BaseStyle = $d.Styles("Fuentedeprrafopredeter")
$ttulo1Car.LinkStyle = $ttulo1
$ttulo1Car.Priority = 9

$ttulo1Car.Font.Name = "Times New Roman"
$ttulo1Car.Font.Bold = $true
$ttulo1Car.Font.Size = 14
$ttulo1Car.Font.SizeBi = 16

Write-Output "done"
